$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").Copy()
$ws.Range("A47").PasteSpecial(-4122)
$ws.Range("A47").Value = 44165
$ws.Range("B47").Value = 5

$ws.Range("B48").Select()
